$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 910472.4  # H2: update (was 1126.4)
$ws.Cells.Item(2, 9).Value = 1233  # I2: update (was 674.875)
$ws.Cells.Item(2, 10).Value = 1251437.1  # J2: update (was 1642.4286)
$ws.Cells.Item(2, 11).Value = 1233  # K2: update (was 674.875)
$ws.Cells.Item(2, 12).Value = 1251437.1  # L2: update (was 1642.4286)
$ws.Cells.Item(2, 13).Value = -1120  # M2: update (was -561.875)
$ws.Cells.Item(2, 14).Value = -1251663.1  # N2: update (was -1868.4286)
$ws.Cells.Item(6, 8).Value = 164.125  # H6: update (was 195.5)
$ws.Cells.Item(6, 9).Value = 182.6  # I6: update (was 195.5)
$ws.Cells.Item(6, 10).Value = 133.33333  # J6: update (was 0)
$ws.Cells.Item(6, 11).Value = 547.8  # K6: update (was 586.5)
$ws.Cells.Item(6, 12).Value = 399.99999  # L6: update (was 0)
$ws.Cells.Item(6, 13).Value = -435.8  # M6: update (was -474.5)
$ws.Cells.Item(6, 14).Value = -623.99999  # N6: add (was None)
$ws.Cells.Item(8, 8).Value = 4765  # H8: update (was 7400)
$ws.Cells.Item(8, 9).Value = 4765  # I8: update (was 7400)
$ws.Cells.Item(8, 11).Value = 14295  # K8: update (was 22200)
$ws.Cells.Item(8, 13).Value = -14156  # M8: update (was -22061)
$ws.Cells.Item(19, 8).Value = 6668412.5  # H19: update (was 4546775.5)
$ws.Cells.Item(19, 9).Value = 1299.6  # I19: update (was 869.1111)
$ws.Cells.Item(19, 10).Value = 10001969  # J19: update (was 7693941.5)
$ws.Cells.Item(19, 11).Value = 1299.6  # K19: update (was 869.1111)
$ws.Cells.Item(19, 12).Value = 10001969  # L19: update (was 7693941.5)
$ws.Cells.Item(19, 13).Value = -1124.6  # M19: update (was -694.1111)
$ws.Cells.Item(19, 14).Value = -10002319  # N19: update (was -7694291.5)
$ws.Cells.Item(21, 8).Value = 0  # H21: update (was 46999)
$ws.Cells.Item(21, 9).Value = 0  # I21: update (was 46999)
$ws.Cells.Item(21, 11).Value = 0  # K21: update (was 46999)
$ws.Cells.Item(21, 13).Value = ""  # M21: remove (was -46531)
$ws.Cells.Item(23, 8).Value = 0  # H23: update (was 46999)
$ws.Cells.Item(23, 9).Value = 0  # I23: update (was 46999)
$ws.Cells.Item(23, 11).Value = 0  # K23: update (was 46999)
$ws.Cells.Item(23, 13).Value = ""  # M23: remove (was -46765)
$ws.Cells.Item(38, 8).Value = 1304.9375  # H38: update (was 1765.5294)
$ws.Cells.Item(38, 9).Value = 134.28572  # I38: update (was 139.76923)
$ws.Cells.Item(38, 10).Value = 9499.5  # J38: update (was 7049.25)
$ws.Cells.Item(38, 11).Value = 402.85716  # K38: update (was 419.30769)
$ws.Cells.Item(38, 12).Value = 28498.5  # L38: update (was 21147.75)
$ws.Cells.Item(38, 13).Value = -30.85716000000002  # M38: update (was -47.30768999999998)
$ws.Cells.Item(38, 14).Value = -29242.5  # N38: update (was -21891.75)
$ws.Cells.Item(43, 8).Value = 3392.182  # H43: update (was 3545.4)
$ws.Cells.Item(43, 10).Value = 3058  # J43: update (was 3229.1428)
$ws.Cells.Item(43, 12).Value = 3058  # L43: update (was 3229.1428)
$ws.Cells.Item(43, 14).Value = -3196  # N43: update (was -3367.1428)
$ws.Cells.Item(51, 8).Value = 7712.5  # H51: update (was 7579.7617)
$ws.Cells.Item(51, 10).Value = 5736.8423  # J51: update (was 5472.222)
$ws.Cells.Item(51, 12).Value = 5736.8423  # L51: update (was 5472.222)
$ws.Cells.Item(51, 14).Value = -6704.8423  # N51: update (was -6440.222)
$ws.Cells.Item(58, 8).Value = 4222.1665  # H58: update (was 4999.4)
$ws.Cells.Item(58, 9).Value = 444.33334  # I58: update (was 498.5)
$ws.Cells.Item(58, 11).Value = 1333.00002  # K58: update (was 1495.5)
$ws.Cells.Item(58, 13).Value = -1183.00002  # M58: update (was -1345.5)
$ws.Cells.Item(100, 8).Value = 5183.4736  # H100: update (was 4974.1)
$ws.Cells.Item(100, 10).Value = 6830.4  # J100: update (was 6300)
$ws.Cells.Item(100, 12).Value = 6830.4  # L100: update (was 6300)
$ws.Cells.Item(100, 14).Value = -7912.4  # N100: update (was -7382)
$ws.Cells.Item(132, 8).Value = 4090.3333  # H132: update (was 4616.303)
$ws.Cells.Item(132, 9).Value = 1939.2222  # I132: update (was 2060.2083)
$ws.Cells.Item(132, 10).Value = 10543.667  # J132: update (was 11432.556)
$ws.Cells.Item(132, 11).Value = 5817.6666  # K132: update (was 6180.624899999999)
$ws.Cells.Item(132, 12).Value = 31631.001  # L132: update (was 34297.66800000001)
$ws.Cells.Item(132, 13).Value = -3287.6666  # M132: update (was -3650.624899999999)
$ws.Cells.Item(132, 14).Value = -36691.001  # N132: update (was -39357.66800000001)
$ws.Cells.Item(137, 8).Value = 1625651.6  # H137: update (was 4333170.5)
$ws.Cells.Item(137, 9).Value = 2073.6  # I137: update (was 5000)
$ws.Cells.Item(137, 10).Value = 4331615  # J137: update (was 6497256)
$ws.Cells.Item(137, 11).Value = 6220.799999999999  # K137: update (was 15000)
$ws.Cells.Item(137, 12).Value = 12994845  # L137: update (was 19491768)
$ws.Cells.Item(137, 13).Value = -3670.799999999999  # M137: update (was -12450)
$ws.Cells.Item(137, 14).Value = -12999945  # N137: update (was -19496868)

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 2662.0833  # H45: update (was 2549.7693)
$ws.Cells.Item(45, 9).Value = 1355.4445  # I45: update (was 1340.1)
$ws.Cells.Item(45, 11).Value = 1355.4445  # K45: update (was 1340.1)
$ws.Cells.Item(45, 13).Value = -978.4445000000001  # M45: update (was -963.0999999999999)
$ws.Cells.Item(61, 8).Value = 23336632  # H61: update (was 19094064)
$ws.Cells.Item(61, 9).Value = 40004580  # I61: update (was 33337484)
$ws.Cells.Item(61, 10).Value = 2501699.8  # J61: update (was 2001959.6)
$ws.Cells.Item(61, 11).Value = 40004580  # K61: update (was 33337484)
$ws.Cells.Item(61, 12).Value = 2501699.8  # L61: update (was 2001959.6)
$ws.Cells.Item(61, 13).Value = -40004368  # M61: update (was -33337272)
$ws.Cells.Item(61, 14).Value = -2502123.8  # N61: update (was -2002383.6)
$ws.Cells.Item(74, 8).Value = 928107.0600000001  # H74: update (was 808411.5)
$ws.Cells.Item(74, 9).Value = 1191875.4  # I74: update (was 1042948.06)
$ws.Cells.Item(74, 10).Value = 4918  # J74: update (was 4286.143)
$ws.Cells.Item(74, 11).Value = 1191875.4  # K74: update (was 1042948.06)
$ws.Cells.Item(74, 12).Value = 4918  # L74: update (was 4286.143)
$ws.Cells.Item(74, 13).Value = -1191001.4  # M74: update (was -1042074.06)
$ws.Cells.Item(74, 14).Value = -6666  # N74: update (was -6034.143)
$ws.Cells.Item(77, 8).Value = 928107.0600000001  # H77: update (was 808411.5)
$ws.Cells.Item(77, 9).Value = 1191875.4  # I77: update (was 1042948.06)
$ws.Cells.Item(77, 10).Value = 4918  # J77: update (was 4286.143)
$ws.Cells.Item(77, 11).Value = 5959377  # K77: update (was 5214740.300000001)
$ws.Cells.Item(77, 12).Value = 24590  # L77: update (was 21430.715)
$ws.Cells.Item(77, 13).Value = -5955009  # M77: update (was -5210372.300000001)
$ws.Cells.Item(77, 14).Value = -33326  # N77: update (was -30166.715)
$ws.Cells.Item(107, 8).Value = 72614  # H107: update (was 0)
$ws.Cells.Item(107, 10).Value = 72614  # J107: update (was 0)
$ws.Cells.Item(107, 12).Value = 72614  # L107: update (was 0)
$ws.Cells.Item(107, 14).Value = -80294  # N107: add (was None)
$ws.Cells.Item(124, 8).Value = 38151.668  # H124: update (was 34613.5)
$ws.Cells.Item(124, 10).Value = 38151.668  # J124: update (was 34613.5)
$ws.Cells.Item(124, 12).Value = 38151.668  # L124: update (was 34613.5)
$ws.Cells.Item(124, 14).Value = -47971.668  # N124: update (was -44433.5)
$ws.Cells.Item(136, 8).Value = 23336632  # H136: update (was 19094064)
$ws.Cells.Item(136, 9).Value = 40004580  # I136: update (was 33337484)
$ws.Cells.Item(136, 10).Value = 2501699.8  # J136: update (was 2001959.6)
$ws.Cells.Item(136, 11).Value = 120013740  # K136: update (was 100012452)
$ws.Cells.Item(136, 12).Value = 7505099.399999999  # L136: update (was 6005878.800000001)
$ws.Cells.Item(136, 13).Value = -120011190  # M136: update (was -100009902)
$ws.Cells.Item(136, 14).Value = -7510199.399999999  # N136: update (was -6010978.800000001)

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 2165  # H94: update (was 2241.9143)
$ws.Cells.Item(94, 9).Value = 1773.5714  # I94: update (was 1885.4615)
$ws.Cells.Item(94, 10).Value = 3161.3635  # J94: update (was 3271.6667)
$ws.Cells.Item(94, 11).Value = 1773.5714  # K94: update (was 1885.4615)
$ws.Cells.Item(94, 12).Value = 3161.3635  # L94: update (was 3271.6667)
$ws.Cells.Item(94, 13).Value = -1322.5714  # M94: update (was -1434.4615)
$ws.Cells.Item(94, 14).Value = -4063.3635  # N94: update (was -4173.6667)
$ws.Cells.Item(134, 8).Value = 11114101  # H134: update (was 12502864)
$ws.Cells.Item(134, 10).Value = 25003000  # J134: update (was 33336000)
$ws.Cells.Item(134, 12).Value = 75009000  # L134: update (was 100008000)
$ws.Cells.Item(134, 14).Value = -75014070  # N134: update (was -100013070)

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(107, 8).Value = 2670.0278  # H107: update (was 2643.1353)
$ws.Cells.Item(107, 9).Value = 2525.577  # I107: update (was 2456.6785)
$ws.Cells.Item(107, 10).Value = 3045.6  # J107: update (was 3223.2222)
$ws.Cells.Item(107, 11).Value = 2525.577  # K107: update (was 2456.6785)
$ws.Cells.Item(107, 12).Value = 3045.6  # L107: update (was 3223.2222)
$ws.Cells.Item(107, 13).Value = -605.5770000000002  # M107: update (was -536.6785)
$ws.Cells.Item(107, 14).Value = -6885.6  # N107: update (was -7063.2222)
$ws.Cells.Item(134, 8).Value = 2206.2144  # H134: update (was 2325.6667)
$ws.Cells.Item(134, 9).Value = 2278.0908  # I134: update (was 2453.3333)
$ws.Cells.Item(134, 11).Value = 6834.2724  # K134: update (was 7359.999899999999)
$ws.Cells.Item(134, 13).Value = -4299.2724  # M134: update (was -4824.999899999999)

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(17, 8).Value = 6447.5713  # H17: update (was 8322)
$ws.Cells.Item(17, 10).Value = 8946.6  # J17: update (was 12383)
$ws.Cells.Item(17, 12).Value = 26839.8  # L17: update (was 37149)
$ws.Cells.Item(17, 14).Value = -27177.8  # N17: update (was -37487)
$ws.Cells.Item(38, 8).Value = 616.5  # H38: update (was 701.4)
$ws.Cells.Item(38, 10).Value = 920.75  # J38: update (was 1163.6666)
$ws.Cells.Item(38, 12).Value = 2762.25  # L38: update (was 3490.9998)
$ws.Cells.Item(38, 14).Value = -3456.25  # N38: update (was -4184.9998)
$ws.Cells.Item(92, 8).Value = 426.4  # H92: update (was 55555910)
$ws.Cells.Item(92, 9).Value = 0  # I92: update (was 333333340)
$ws.Cells.Item(92, 11).Value = 0  # K92: update (was 1000000020)
$ws.Cells.Item(92, 13).Value = ""  # M92: remove (was -999998772)
$ws.Cells.Item(112, 8).Value = 18335.8  # H112: update (was 21388)
$ws.Cells.Item(112, 9).Value = 10013  # I112: update (was 14995)
$ws.Cells.Item(112, 10).Value = 20416.5  # J112: update (was 22666.6)
$ws.Cells.Item(112, 11).Value = 30039  # K112: update (was 44985)
$ws.Cells.Item(112, 12).Value = 61249.5  # L112: update (was 67999.79999999999)
$ws.Cells.Item(112, 13).Value = -28931  # M112: update (was -43877)
$ws.Cells.Item(112, 14).Value = -63465.5  # N112: update (was -70215.79999999999)
$ws.Cells.Item(137, 8).Value = 5327.7144  # H137: update (was 3930.8096)
$ws.Cells.Item(137, 9).Value = 2771.25  # I137: update (was 2234.111)
$ws.Cells.Item(137, 10).Value = 20666.5  # J137: update (was 14111)
$ws.Cells.Item(137, 11).Value = 8313.75  # K137: update (was 6702.333)
$ws.Cells.Item(137, 12).Value = 61999.5  # L137: update (was 42333)
$ws.Cells.Item(137, 13).Value = -3213.75  # M137: update (was -1602.333)
$ws.Cells.Item(137, 14).Value = -72199.5  # N137: update (was -52533)
$ws.Cells.Item(141, 8).Value = 8560.286  # H141: update (was 7940.125)
$ws.Cells.Item(141, 9).Value = 4431.5  # I141: update (was 4312.5713)
$ws.Cells.Item(141, 11).Value = 13294.5  # K141: update (was 12937.7139)
$ws.Cells.Item(141, 13).Value = -8114.5  # M141: update (was -7757.713899999999)

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 801.5714  # H97: update (was 856.1429000000001)
$ws.Cells.Item(97, 9).Value = 691.5  # I97: update (was 755.6667)
$ws.Cells.Item(97, 10).Value = 948.3333  # J97: update (was 931.5)
$ws.Cells.Item(97, 11).Value = 691.5  # K97: update (was 755.6667)
$ws.Cells.Item(97, 12).Value = 948.3333  # L97: update (was 931.5)
$ws.Cells.Item(97, 13).Value = -195.5  # M97: update (was -259.6667)
$ws.Cells.Item(97, 14).Value = -1940.3333  # N97: update (was -1923.5)
$ws.Cells.Item(102, 8).Value = 2081.92  # H102: update (was 2145.28)
$ws.Cells.Item(102, 9).Value = 2072.875  # I102: update (was 2138.875)
$ws.Cells.Item(102, 11).Value = 2072.875  # K102: update (was 2138.875)
$ws.Cells.Item(102, 13).Value = -450.875  # M102: update (was -516.875)
$ws.Cells.Item(122, 8).Value = 4634.1113  # H122: update (was 4397)
$ws.Cells.Item(122, 9).Value = 5744.5  # I122: update (was 4995.4)
$ws.Cells.Item(122, 10).Value = 3745.8  # J122: update (was 3798.6)
$ws.Cells.Item(122, 11).Value = 17233.5  # K122: update (was 14986.2)
$ws.Cells.Item(122, 12).Value = 11237.4  # L122: update (was 11395.8)
$ws.Cells.Item(122, 13).Value = -14783.5  # M122: update (was -12536.2)
$ws.Cells.Item(122, 14).Value = -16137.4  # N122: update (was -16295.8)
$ws.Cells.Item(123, 8).Value = 91713.57000000001  # H123: update (was 99999.336)
$ws.Cells.Item(123, 10).Value = 91713.57000000001  # J123: update (was 99999.336)
$ws.Cells.Item(123, 12).Value = 91713.57000000001  # L123: update (was 99999.336)
$ws.Cells.Item(123, 14).Value = -96613.57000000001  # N123: update (was -104899.336)
$ws.Cells.Item(126, 8).Value = 1512  # H126: update (was 1510.5)
$ws.Cells.Item(126, 10).Value = 0  # J126: update (was 1500)
$ws.Cells.Item(126, 12).Value = 0  # L126: update (was 4500)
$ws.Cells.Item(126, 14).Value = ""  # N126: remove (was -9440)

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 3328.818  # H16: update (was 3186.0435)
$ws.Cells.Item(16, 9).Value = 1805.5  # I16: update (was 1708.0588)
$ws.Cells.Item(16, 10).Value = 7391  # J16: update (was 7373.6665)
$ws.Cells.Item(16, 11).Value = 1805.5  # K16: update (was 1708.0588)
$ws.Cells.Item(16, 12).Value = 7391  # L16: update (was 7373.6665)
$ws.Cells.Item(16, 13).Value = -1635.5  # M16: update (was -1538.0588)
$ws.Cells.Item(16, 14).Value = -7731  # N16: update (was -7713.6665)
$ws.Cells.Item(93, 8).Value = 11122911  # H93: update (was 3090299.5)
$ws.Cells.Item(93, 9).Value = 0  # I93: update (was 894.4545000000001)
$ws.Cells.Item(93, 10).Value = 11122911  # J93: update (was 7945079)
$ws.Cells.Item(93, 11).Value = 0  # K93: update (was 894.4545000000001)
$ws.Cells.Item(93, 12).Value = 11122911  # L93: update (was 7945079)
$ws.Cells.Item(93, 13).Value = ""  # M93: remove (was 353.5454999999999)
$ws.Cells.Item(93, 14).Value = -11125407  # N93: update (was -7947575)
$ws.Cells.Item(100, 8).Value = 27783038  # H100: update (was 31255292)
$ws.Cells.Item(100, 10).Value = 62505748  # J100: update (was 83339336)
$ws.Cells.Item(100, 12).Value = 62505748  # L100: update (was 83339336)
$ws.Cells.Item(100, 14).Value = -62506830  # N100: update (was -83340418)

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(103, 8).Value = 36911.832  # H103: update (was 40000.2)
$ws.Cells.Item(103, 10).Value = 36911.832  # J103: update (was 40000.2)
$ws.Cells.Item(103, 12).Value = 36911.832  # L103: update (was 40000.2)
$ws.Cells.Item(103, 14).Value = -39255.832  # N103: update (was -42344.2)
